$wb = $excel.ActiveWorkbook

# Rename category labels: "light goods" -> "van", "heavy goods" -> "lorry"
# These labels live in column A, row 3 and row 4 on every data sheet.
foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A3").Value2 -eq "light goods") {
        $ws.Range("A3").Value = "van"
    }
    if ($ws.Range("A4").Value2 -eq "heavy goods") {
        $ws.Range("A4").Value = "lorry"
    }
}


# --- Sheet: mean ---
$ws = $wb.Worksheets.Item("mean")
$ws.Range("B2").Value = 3.6941528940735493
$ws.Range("C2").Value = 3.116195746747314
$ws.Range("D2").Value = 3.989104174745653
$ws.Range("E2").Value = 2.4034183688700517
$ws.Range("F2").Value = 3.2547603362365223
$ws.Range("B3").Value = 3.081163271587731
$ws.Range("C3").Value = 2.4382992006160045
$ws.Range("D3").Value = 2.862070987505259
$ws.Range("E3").Value = 1.9119959253232333
$ws.Range("F3").Value = 2.5955470064007256
$ws.Range("B4").Value = 15.873055846276467
$ws.Range("C4").Value = 18.78433824470811
$ws.Range("D4").Value = 19.77666465776771
$ws.Range("E4").Value = 19.02291067702824
$ws.Range("F4").Value = 17.10531817759326
$ws.Range("B5").Value = 9.236023440704061
$ws.Range("C5").Value = 11.574977021548094
$ws.Range("D5").Value = 8.257702563934133
$ws.Range("E5").Value = 4.101707673923155
$ws.Range("F5").Value = 7.689525045009852
$ws.Range("B6").Value = 22.830163553082294
$ws.Range("C6").Value = 25.315402881694617
$ws.Range("D6").Value = 18.783759866192973
$ws.Range("E6").Value = 13.545189985488175
$ws.Range("F6").Value = 19.260592274824077
$ws.Range("B7").Value = 5.086910163266265
$ws.Range("C7").Value = 2.0854692494558797
$ws.Range("D7").Value = 1.3080194394595497
$ws.Range("E7").Value = 0.6405785004784111
$ws.Range("F7").Value = 1.1342981864809742

# --- Sheet: median ---
$ws = $wb.Worksheets.Item("median")
$ws.Range("B2").Value = 3.6937273529915524
$ws.Range("C2").Value = 3.116246408791473
$ws.Range("D2").Value = 3.992409031793118
$ws.Range("E2").Value = 2.4030744172215437
$ws.Range("F2").Value = 3.254222463621489
$ws.Range("B3").Value = 3.083886547937352
$ws.Range("C3").Value = 2.440970871682321
$ws.Range("D3").Value = 2.8551700803774684
$ws.Range("E3").Value = 1.9091716431730297
$ws.Range("F3").Value = 2.592433958499204
$ws.Range("B4").Value = 15.868064863422077
$ws.Range("C4").Value = 18.76640492276497
$ws.Range("D4").Value = 19.737857515596275
$ws.Range("E4").Value = 18.998366123902084
$ws.Range("F4").Value = 17.107164662443154
$ws.Range("B5").Value = 9.17800066884514
$ws.Range("C5").Value = 11.535931985043705
$ws.Range("D5").Value = 8.210616933313968
$ws.Range("E5").Value = 4.081454768948477
$ws.Range("F5").Value = 7.687759895905983
$ws.Range("B6").Value = 22.710845304995324
$ws.Range("C6").Value = 25.250585783853687
$ws.Range("D6").Value = 18.68087757318349
$ws.Range("E6").Value = 13.54356827092101
$ws.Range("F6").Value = 19.239353402432847
$ws.Range("B7").Value = 4.7358412707244675
$ws.Range("C7").Value = 2.0450897995898507
$ws.Range("D7").Value = 1.2795972939776212
$ws.Range("E7").Value = 0.630545229843144
$ws.Range("F7").Value = 1.128329419942393

# --- Sheet: lower 5 ---
$ws = $wb.Worksheets.Item("lower 5")
$ws.Range("B2").Value = 3.593958656889468
$ws.Range("C2").Value = 3.0097367866632094
$ws.Range("D2").Value = 3.837162613946804
$ws.Range("E2").Value = 2.317900416710358
$ws.Range("F2").Value = 3.201686370675815
$ws.Range("B3").Value = 2.856885408002161
$ws.Range("C3").Value = 2.1909200378455176
$ws.Range("D3").Value = 2.6013053831465234
$ws.Range("E3").Value = 1.723418262259683
$ws.Range("F3").Value = 2.4841614689420344
$ws.Range("B4").Value = 15.17074088803298
$ws.Range("C4").Value = 17.426732598547304
$ws.Range("D4").Value = 17.904038662690066
$ws.Range("E4").Value = 16.980688977830454
$ws.Range("F4").Value = 16.545883161703998
$ws.Range("B5").Value = 7.833564835511852
$ws.Range("C5").Value = 9.768062738474839
$ws.Range("D5").Value = 6.608737960317089
$ws.Range("E5").Value = 3.309934748087556
$ws.Range("F5").Value = 7.02016497189735
$ws.Range("B6").Value = 20.018245957901065
$ws.Range("C6").Value = 23.034191865079663
$ws.Range("D6").Value = 15.631866506366196
$ws.Range("E6").Value = 12.022685780156541
$ws.Range("F6").Value = 18.155374890454812
$ws.Range("B7").Value = 2.4394054483147243
$ws.Range("C7").Value = 1.3146331230914994
$ws.Range("D7").Value = 0.8058288173212159
$ws.Range("E7").Value = 0.4193358606355696
$ws.Range("F7").Value = 0.8915037414127244

# --- Sheet: upper 95 ---
$ws = $wb.Worksheets.Item("upper 95")
$ws.Range("B2").Value = 3.801126202605852
$ws.Range("C2").Value = 3.2335543679234746
$ws.Range("D2").Value = 4.141486241880108
$ws.Range("E2").Value = 2.490032227956992
$ws.Range("F2").Value = 3.3084263161844074
$ws.Range("B3").Value = 3.2980202643835113
$ws.Range("C3").Value = 2.700545776023522
$ws.Range("D3").Value = 3.1408338873226738
$ws.Range("E3").Value = 2.1132693217422336
$ws.Range("F3").Value = 2.710011973914286
$ws.Range("B4").Value = 16.569697432350722
$ws.Range("C4").Value = 20.29357807515842
$ws.Range("D4").Value = 21.776798180977327
$ws.Range("E4").Value = 20.987738259530566
$ws.Range("F4").Value = 17.69597721555329
$ws.Range("B5").Value = 10.85615301409685
$ws.Range("C5").Value = 13.44612101417834
$ws.Range("D5").Value = 10.100616934402128
$ws.Range("E5").Value = 4.97629352569548
$ws.Range("F5").Value = 8.406592746650333
$ws.Range("B6").Value = 25.78934047249593
$ws.Range("C6").Value = 27.806236893616507
$ws.Range("D6").Value = 22.326689047482002
$ws.Range("E6").Value = 15.127470666100837
$ws.Range("F6").Value = 20.413120054707615
$ws.Range("B7").Value = 8.78428830368854
$ws.Range("C7").Value = 3.0860266315067113
$ws.Range("D7").Value = 1.9342917255887289
$ws.Range("E7").Value = 0.8948563599626743
$ws.Range("F7").Value = 1.3995784706992287

